$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Subject Code column (B2:B4) from "DAM-101" to "DME-105"
$ws.Range("B2").Value = "DME-105"
$ws.Range("B3").Value = "DME-105"
$ws.Range("B4").Value = "DME-105"

# Reflect the last active cell selection as seen in the saved file
$ws.Range("B4").Select()
